{"js": "// Update the \"12-month % Survival\" column (last column) for several rows\n// that currently show \"NA\" with their now-known survival percentages, and\n// refresh the \"Totals and Weighted Averages\" row's weighted-average cell.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst values = table.values;\n\n// Map Site ID (column index 2) -> new 12-month % Survival value (last column).\nconst updatesBySiteId = {\n  \"IC_U_3_IM_1\": \"74\",\n  \"IC_C_3_IM_1\": \"97\",\n  \"IC_Z_3_IM_1\": \"100\",\n  \"IC_C_3_M_2\": \"84.2\",\n  \"IC_U_3_M_2\": \"59.2\",\n  \"IC_Z_3_M_2\": \"57.5\",\n  \"IC_U_3_IM_2\": \"64.83\",\n  \"IC_C_3_IM_2\": \"80.67\",\n  \"IC_Z_3_IM_2\": \"84\",\n};\n\nconst lastCol = values[0].length - 1;\n\nfor (let r = 0; r < values.length; r++) {\n  const siteId = values[r][2];\n  if (Object.prototype.hasOwnProperty.call(updatesBySiteId, siteId)) {\n    table.getCell(r, lastCol).value = updatesBySiteId[siteId];\n  }\n}\n\n// Totals / weighted-average row is the last row in the table.\nconst totalsRow = values.length - 1;\ntable.getCell(totalsRow, lastCol).value = \"70.48\";\n\nawait context.sync();\n", "ps1": "# Update the \"12-month % Survival\" column (last column) for several rows\n# that currently show \"NA\" with their now-known survival percentages, and\n# refresh the \"Totals and Weighted Averages\" row's weighted-average cell.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$lastCol = $t.Columns.Count\n\n$updatesBySiteId = @{\n    \"IC_U_3_IM_1\" = \"74\";\n    \"IC_C_3_IM_1\" = \"97\";\n    \"IC_Z_3_IM_1\" = \"100\";\n    \"IC_C_3_M_2\"  = \"84.2\";\n    \"IC_U_3_M_2\"  = \"59.2\";\n    \"IC_Z_3_M_2\"  = \"57.5\";\n    \"IC_U_3_IM_2\" = \"64.83\";\n    \"IC_C_3_IM_2\" = \"80.67\";\n    \"IC_Z_3_IM_2\" = \"84\";\n}\n\n$rowCount = $t.Rows.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    # Cell.Range.Text carries a trailing cell-mark (CR + BEL); strip it before\n    # comparing so the Site ID matches cleanly.\n    $siteId = $t.Cell($r, 3).Range.Text.TrimEnd([char]13, [char]7)\n    if ($updatesBySiteId.ContainsKey($siteId)) {\n        $t.Cell($r, $lastCol).Range.Text = $updatesBySiteId[$siteId]\n    }\n}\n\n# Totals / weighted-average row is the last row in the table.\n$t.Cell($rowCount, $lastCol).Range.Text = \"70.48\"\n"}
